$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 27
$url = "https://www.360dx.com/regulatory-news-fda-approvals/roche-abbott-hologic-cepheid-others-gain-510k-clearances-february"
$keyword = "digital pathology"
$title = "Roche, Abbott, Hologic, Cepheid, Others Gain 510(k) Clearances in February"

$ws.Hyperlinks.Add($ws.Range("A$newRow"), $url) | Out-Null
$ws.Range("A$newRow").Style = "Hyperlink"
$ws.Range("B$newRow").Value = $keyword
$ws.Range("C$newRow").Value = $title
